# Add a new HARWIN connector row to the component list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 5

$ws.Range("A$newRow").Value = "Receptacle"
$ws.Range("B$newRow").Value = "2x8"
$ws.Range("C$newRow").Value = "SMD"
$ws.Range("D$newRow").Value = "HARWIN"
$ws.Range("F$newRow").Value = "M20-7870842"

$linkCell = $ws.Range("I$newRow")
$linkCell.Value = "http://www.harwin.com/search/M20-7870842?ProductSearch=True"
$ws.Hyperlinks.Add($linkCell, "http://www.harwin.com/search/M20-7870842?ProductSearch=True")
$linkCell.Style = "Hyperlink"

$ws.Range("A6").Select()
